$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.809.16"
Set-TextValue "E2" "  -1.51%  "
Set-TextValue "D3" "1.891.71"
Set-TextValue "E3" "  -1.21%  "
Set-TextValue "D4" "1.0000"
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "0.7637"
Set-TextValue "E5" "  +2.77%  "
Set-TextValue "D6" "238.96"
Set-TextValue "E6" "  -2.28%  "
Set-TextValue "D7" "0.9993"
Set-TextValue "E7" "  +0.01%  "
Set-TextValue "D8" "0.3042"
Set-TextValue "E8" "  -3.58%  "
Set-TextValue "D9" "25.32"
Set-TextValue "E9" "  -7.28%  "
Set-TextValue "D10" "0.06807"
Set-TextValue "E10" "  -2.72%  "
Set-TextValue "D11" "0.07989"
Set-TextValue "E11" "  +0.27%  "
Set-TextValue "D12" "0.7469"
Set-TextValue "E12" "  -4.68%  "
Set-TextValue "D13" "1.888.05"
Set-TextValue "E13" "  -1.65%  "
Set-TextValue "D14" "5.194"
Set-TextValue "E14" "  -1.92%  "
Set-TextValue "D15" "90.91"
Set-TextValue "E15" "  -1.16%  "
Set-TextValue "D16" "29.808.33"
Set-TextValue "E16" "  -1.85%  "
Set-TextValue "E17" "  -3.44%  "
Set-TextValue "D18" "5.947"
Set-TextValue "E18" "  +1.75%  "
Set-TextValue "D19" "0.000007663"
Set-TextValue "E19" "  -2.25%  "
Set-TextValue "D20" "234.42"
Set-TextValue "E20" "  -4.84%  "
Set-TextValue "D21" "0.9994"
Set-TextValue "E21" "  -0.04%  "
Set-TextValue "D22" "2.133.92"
Set-TextValue "E22" "  -4.34%  "
Set-TextValue "D23" "1.000"
Set-TextValue "E23" "  +0.16%  "
Set-TextValue "D24" "6.924"
Set-TextValue "E24" "  +3.68%  "
Set-TextValue "D25" "9.221"
Set-TextValue "E25" "  -2.51%  "
Set-TextValue "D26" "165.33"
Set-TextValue "E26" "  +0.18%  "
Set-TextValue "D27" "18.69"
Set-TextValue "E27" "  -2.21%  "
Set-TextValue "D28" "0.1307"
Set-TextValue "E28" "  +2.30%  "
Set-TextValue "D29" "2.040"
Set-TextValue "E29" "  -4.84%  "
Set-TextValue "D30" "1.341"
Set-TextValue "E30" "  -1.29%  "
Set-TextValue "D31" "1.508"
Set-TextValue "E31" "  -2.68%  "
Set-TextValue "D32" "4.270"
Set-TextValue "E32" "  -1.58%  "
Set-TextValue "D33" "4.014"
Set-TextValue "E33" "  -1.94%  "
Set-TextValue "D34" "0.05363"
Set-TextValue "E34" "  +2.26%  "
Set-TextValue "D35" "1.245"
Set-TextValue "E35" "  -4.64%  "
Set-TextValue "D36" "0.7261"
Set-TextValue "E36" "  -3.77%  "
Set-TextValue "D37" "2.710"
Set-TextValue "E37" "  -1.73%  "
Set-TextValue "D38" "0.01922"
Set-TextValue "E38" "  -1.31%  "
Set-TextValue "D39" "2.772"
Set-TextValue "E39" "  -0.49%  "
Set-TextValue "D40" "6.172"
Set-TextValue "E40" "  -3.86%  "
Set-TextValue "D41" "0.4397"
Set-TextValue "E41" "  -2.52%  "
Set-TextValue "D42" "72.03"
Set-TextValue "E42" "  -5.59%  "
Set-TextValue "D43" "1.909"
Set-TextValue "E43" "  -2.76%  "
Set-TextValue "D44" "0.9992"
Set-TextValue "E44" "  +0.06%  "
Set-TextValue "E45" "  -1.21%  "
Set-TextValue "D46" "100.93"
Set-TextValue "E46" "  -0.52%  "
Set-TextValue "D47" "7.570"
Set-TextValue "E47" "  -2.63%  "
Set-TextValue "D48" "9.790"
Set-TextValue "E48" "  -1.21%  "
Set-TextValue "D49" "2.035.93"
Set-TextValue "E49" "  -4.02%  "
Set-TextValue "D50" "36.16"
Set-TextValue "E50" "  -2.78%  "
Set-TextValue "D51" "923.58"
Set-TextValue "E51" "  -2.51%  "
